$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fuel salt (ThF4) now has a real price instead of "n/a", plus a unit-price
# reference value in column F.
$ws.Range("B4").Value = 152700
$ws.Range("F4").Value = 26017

# New reference row showing the UF6 price basis used for B4/F4 above.
$ws.Range("E6").Value = "UF6"
$ws.Range("F6").Value = "152700/m^3"

# Fuel-salt cost buildup now includes the UF4 contribution (previously
# blank / driving a bogus #VALUE! side calculation instead).
$ws.Range("B11").Formula = "=B4*A11"
$ws.Range("B18").Formula = "=B4*A18"

# Remove the old broken side-table (columns I/J/K, rows 9-14) that tried to
# multiply text ("n/a"/"total"/"$/m^3"/"$/ft^3") by numbers and produced
# #VALUE! errors. Clear() wipes both contents and the cell's style so the
# emptied cells (and, where a whole row empties out, the row itself) drop
# out of the saved worksheet entirely.
$ws.Range("I9").Clear()
$ws.Range("I10:J10").Clear()
$ws.Range("I11:J11").Clear()
$ws.Range("I12:J12").Clear()
$ws.Range("I13:K13").Clear()
$ws.Range("J14:K14").Clear()

# Match the author's final selection position.
$ws.Range("K18").Select() | Out-Null
